$d = $word.ActiveDocument

# 1. "foeuilles" -> "foeilles"
$d.Content.Find.Execute("toutes les foeuilles, si tu veulx,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "toutes les foeilles, si tu veulx,", 2)

# 2. "doict" -> "doibt"
$d.Content.Find.Execute("qui se doict brusler ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "qui se doibt brusler ", 2)

# 3. "escrevisse" -> "escrevisses"
$d.Content.Find.Execute("escrevisse", $true, $false, $false, $false, $false,
                         $true, 1, $false, "escrevisses", 2)

# 4. "Et quand ils sont une foys recuits, ne les garde gueres"
#    -> split into three runs: "Et quand il" + "z" (no direct color - plain run) +
#       " sont une foys recuits, ne les garde gueres"
$target = $d.Content
$target.Find.Execute("Et quand ils sont une foys recuits, ne les garde gueres", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$sentenceStart = $target.Start

# Build a throwaway "z" at this paragraph's own start: inserting right before the
# first run of a paragraph yields a bare run with no inherited direct formatting
# (no color). We then cut that bare run and paste it into the gap left by removing
# the "s" from "ils", giving the desired 3-way run split without carrying over the
# surrounding run's color.
$para = $d.Paragraphs.First
$paraStart = $target.Paragraphs.First.Range.Start

$stub = $d.Range($paraStart, $paraStart)
$stub.InsertBefore("z")
$stubRange = $d.Range($paraStart, $paraStart + 1)
$stubRange.Cut()

# After the cut, the sentence shifted back by 1 (the stub char was removed from
# before it), so re-find the sentence fresh to get reliable offsets.
$target2 = $d.Content
$target2.Find.Execute("Et quand ils sont une foys recuits, ne les garde gueres", $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)
$start2 = $target2.Start

# "Et quand il" is 11 characters; the "s" of "ils" is the next character.
$sChar = $d.Range($start2 + 11, $start2 + 12)
$sChar.Delete()

$gap = $d.Range($start2 + 11, $start2 + 11)
$gap.Paste()
